$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B3" = 6.461
    "E3" = 16.235
    "A12" = -21.564
    "B14" = 6.082
    "E20" = 16.319
    "E25" = 17.287
    "B26" = 6.720999999999999
    "A27" = -21.507
    "E30" = 16.152
    "B31" = 6.462000000000001
    "A32" = -21.352
    "B35" = 7.342000000000001
    "A36" = -21.112
    "B37" = 7.848000000000002
    "A38" = -20.354
    "E44" = 16.317
    "B45" = 6.061
    "A46" = -21.532
    "E47" = 16.582
    "B52" = 5.359
    "A54" = -21.862
    "A55" = -22.247
    "A56" = -21.816
    "B57" = 5.331999999999999
    "E58" = 16.58
    "A67" = -21.623
    "A69" = -21.656
    "A72" = -21.57
    "E78" = 16.564
    "B81" = 6.845999999999999
    "A83" = -20.15
    "B83" = 7.353000000000002
    "E84" = 16.556
    "A86" = -22.096
    "E89" = 17.289
    "A91" = -21.508
    "E91" = 17.39
    "E92" = 17.262
    "A93" = -21.49
    "E96" = 16.545
    "A99" = -20.692
    "B100" = 5.459000000000001
    "B102" = 7.279999999999999
    "E102" = 16.249
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
